$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.81
$ws.Range("G2").Value = 1.86
$ws.Range("I2").Value = 4.7
$ws.Range("J2").Value = 4.2
$ws.Range("N2").Value = 5.4
$ws.Range("P2").Value = 2.5
$ws.Range("R2").Value = 1.61
$ws.Range("U2").Value = 2.42
$ws.Range("AE2").Value = 1000
$ws.Range("AL2").Value = 27
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 9.199999999999999
$ws.Range("AO2").Value = 38
$ws.Range("F3").Value = 1.43
$ws.Range("G3").Value = 1.44
$ws.Range("J3").Value = 5.2
$ws.Range("K3").Value = 5.3
$ws.Range("L3").Value = 1.32
$ws.Range("N3").Value = 5.5
$ws.Range("P3").Value = 2.46
$ws.Range("Z3").Value = 80
$ws.Range("AA3").Value = 290
$ws.Range("AD3").Value = 32
$ws.Range("AF3").Value = 8.800000000000001
$ws.Range("AH3").Value = 24
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 3.45
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 2.68
$ws.Range("W4").Value = 1.83
$ws.Range("AC4").Value = 8.800000000000001
$ws.Range("F5").Value = 1.5
$ws.Range("G5").Value = 1.64
$ws.Range("H5").Value = 5.2
$ws.Range("N5").Value = 6
$ws.Range("P5").Value = 2.94
$ws.Range("Q5").Value = 1.39
$ws.Range("R5").Value = 1.79
$ws.Range("S5").Value = 1.98
$ws.Range("W5").Value = 2.56
$ws.Range("Y5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 5.4
$ws.Range("AO5").Value = 1000
$ws.Range("I6").Value = 2.66
$ws.Range("L6").Value = 1.24
$ws.Range("V6").Value = 1.6
$ws.Range("AF6").Value = 27
$ws.Range("F7").Value = 1.5
$ws.Range("I7").Value = 9
$ws.Range("Q7").Value = 1.71
$ws.Range("S7").Value = 2.8
$ws.Range("W7").Value = 2.68
$ws.Range("F8").Value = 2.52
$ws.Range("I8").Value = 3.3
$ws.Range("L8").Value = 1.48
$ws.Range("N8").Value = 3.5
$ws.Range("Z8").Value = 20
$ws.Range("AA8").Value = 55
$ws.Range("F9").Value = 2.86
$ws.Range("G9").Value = 2.88
$ws.Range("H9").Value = 2.5
$ws.Range("I9").Value = 2.52
$ws.Range("J9").Value = 3.9
$ws.Range("K9").Value = 3.95
$ws.Range("R9").Value = 1.74
$ws.Range("V9").Value = 1.65
$ws.Range("W9").Value = 1.53
$ws.Range("X9").Value = 26
$ws.Range("Z9").Value = 20
$ws.Range("AC9").Value = 9.800000000000001
$ws.Range("AI9").Value = 26
$ws.Range("P10").Value = 2.54
$ws.Range("R10").Value = 1.61
$ws.Range("S10").Value = 2.58
$ws.Range("T10").Value = 1.95
$ws.Range("AJ10").Value = 340
$ws.Range("AL10").Value = 100
$ws.Range("F11").Value = 1.19
$ws.Range("G11").Value = 1.2
$ws.Range("H11").Value = 21
$ws.Range("J11").Value = 8.6
$ws.Range("K11").Value = 8.800000000000001
$ws.Range("N11").Value = 6.6
$ws.Range("T11").Value = 2.34
$ws.Range("W11").Value = 6
$ws.Range("X11").Value = 32
$ws.Range("Y11").Value = 65
$ws.Range("AC11").Value = 19.5
$ws.Range("AD11").Value = 1000
$ws.Range("AH11").Value = 46
$ws.Range("AJ11").Value = 8.6
$ws.Range("H12").Value = 11
$ws.Range("P12").Value = 3.6
$ws.Range("R12").Value = 2.06
$ws.Range("T12").Value = 1.72
$ws.Range("AB12").Value = 16
$ws.Range("F13").Value = 5.6
$ws.Range("G13").Value = 5.7
$ws.Range("H13").Value = 1.69
$ws.Range("I13").Value = 1.7
$ws.Range("V13").Value = 2.42
$ws.Range("W13").Value = 1.21
$ws.Range("Z13").Value = 10
$ws.Range("AB13").Value = 19.5
$ws.Range("AD13").Value = 9.6
$ws.Range("AF13").Value = 44
$ws.Range("AO13").Value = 9
$ws.Range("G14").Value = 3.15
$ws.Range("I14").Value = 2.46
$ws.Range("P14").Value = 2.3
$ws.Range("T14").Value = 1.62
$ws.Range("V14").Value = 1.68
$ws.Range("F15").Value = 2.42
$ws.Range("G15").Value = 2.86
$ws.Range("N15").Value = 5.1
$ws.Range("P15").Value = 2.4
$ws.Range("S15").Value = 2.34
$ws.Range("T15").Value = 1.5
$ws.Range("AA15").Value = 1000
$ws.Range("AK15").Value = 26
$ws.Range("AL15").Value = 32
$ws.Range("AO15").Value = 16.5
$ws.Range("G16").Value = 2.86
$ws.Range("I16").Value = 3.2
$ws.Range("J16").Value = 3.2
$ws.Range("M16").Value = 1.08
$ws.Range("Q16").Value = 2.02
$ws.Range("T16").Value = 1.75
$ws.Range("X16").Value = 15.5
$ws.Range("Z16").Value = 24
$ws.Range("AE16").Value = 38
